$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M2: cited_by_count 25 -> 26
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "26"

# Rows 4, 5, 6 get cyclically rotated:
#   new row 4 <- old row 6
#   new row 5 <- old row 4
#   new row 6 <- old row 5

$oldRow4 = @(
    "José Ribeiro Ferreira, I. Simões, A. Paulo Coimbra, Manuel Crisóstomo",
    "Department of Electrical Engineering, Superior Institute of Engineering of Coimbra, Coimbra, Portugal; Department of Electrical Engineering, Superior Institute of Engineering of Coimbra, Coimbra, Portugal; ISR - Department of Electrical and Computer Engineering, University of Coimbra, Coimbra, Portugal; ISR - Department of Electrical and Computer Engineering, University of Coimbra, Coimbra, Portugal",
    "https://openalex.org/W4226410730",
    "Human–Machine Interfaces to NX100 Controller for Motoman HP3L Robot in Educational Environment",
    "2022-01-01",
    "Communications in computer and information science",
    "N/A",
    "https://doi.org/10.1007/978-3-031-03884-6_3",
    "N/A",
    "N/A",
    "closed",
    "en",
    "0",
    "2022",
    "NA",
    "https://doi.org/10.1007/978-3-031-03884-6_3",
    "book-chapter"
)

$oldRow5 = @(
    "Tao Liu, José Ribeiro Ferreira",
    "State Key Laboratory of Fluid Power and Mechatronic Systems, School of Mechanical Engineering, Zhejiang University, Hangzhou 310027, China; Electrical Engineering Department, Superior Institute of Engineering of Coimbra, 3030-199 Coimbra, Portugal",
    "https://openalex.org/W4285800324",
    "Editorial for the Special Issue on Physical Diagnosis and Rehabilitation Technologies",
    "2022-07-18",
    "Electronics",
    "Multidisciplinary Digital Publishing Institute",
    "https://doi.org/10.3390/electronics11142247",
    "cc-by",
    "publishedVersion",
    "gold",
    "en",
    "0",
    "2022",
    "NA",
    "https://doi.org/10.3390/electronics11142247",
    "article"
)

$oldRow6 = @(
    "Xiaopeng Huang, Tao Liu, Meimei Han, José Ribeiro Ferreira",
    "State Key Laboratory of Fluid Power and Mechatronic Systems, School of Mechanical Engineering, Zhejiang University,Hangzhou,China,310027; State Key Laboratory of Fluid Power and Mechatronic Systems, School of Mechanical Engineering, Zhejiang University,Hangzhou,China,310027; Zhejiang Fuzhi Science and Technology Innovation Co., Ltd.,Hangzhou,China,310027; Institute of Superior of Engineering of Coimbra,Coimbra,Portugal,3030-199",
    "https://openalex.org/W4378965979",
    "Virtual Model Control for Dynamic Banlance of a Two Wheeled-legged Robot",
    "2022-06-01",
    "N/A",
    "N/A",
    "https://doi.org/10.1109/icosr57188.2022.00016",
    "N/A",
    "N/A",
    "closed",
    "en",
    "0",
    "2022",
    "NA",
    "https://doi.org/10.1109/icosr57188.2022.00016",
    "article"
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

function Set-RowValues($rowNum, $values) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $cell = $ws.Range("$col$rowNum")
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
    }
}

Set-RowValues 4 $oldRow6
Set-RowValues 5 $oldRow4
Set-RowValues 6 $oldRow5
